# Applies the three learning-outcome / exam-info text revisions described
# by the commit. (A fourth, purely cosmetic <w:lastRenderedPageBreak/>
# relocation in the diff is a Word-internal repagination bookkeeping
# artifact with no COM-exposed surface -- not something an editing script
# produces -- so it is intentionally left alone.)

$d = $word.ActiveDocument

# 1) Learning outcome bullet: swap out the old "Resonnere ..." sentence for
#    the new "Forklare og redegjøre ..." phrasing.
$d.Content.Find.Execute(
    "Resonnere til hvordan fremtidige ingeniøroppgaver ved bruk av digitale verktøy og nye arbeidsprosesser",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Forklare og redegjøre for ingeniørens rolle i nye arbeidsprosesser og ved økt bruk av digitale verktøy ",
    2) | Out-Null

# 2) "... kvalifisert til sluttkarakter i emnet." -> "... kvalifisert til
#    vurdering i emnet."
$d.Content.Find.Execute(
    "kvalifisert til sluttkarakter i emnet.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "kvalifisert til vurdering i emnet.",
    2) | Out-Null

# 3) Append the pass/fail grading note right after "Det er ingen eksamen i
#    faget."
$d.Content.Find.Execute(
    "Det er ingen eksamen i faget. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Det er ingen eksamen i faget. Bestått/ikke bestått vil gis etter avsluttende muntlig gruppepresentasjon ",
    2) | Out-Null
